$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (Sending cluster = ECs) with recalculated values ---
# Row 2 (Target cluster = ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.699506666666667
$ws.Range("H2").Value = 14.09852
$ws.Range("I2").Value = 0.9660495246229048
$ws.Range("J2").Value = 0.9660495246229047
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.82741333333333
$ws.Range("N2").Value = 95.48223999999999
$ws.Range("O2").Value = 0.114390792932228
$ws.Range("P2").Value = 0.114390792932228
$ws.Range("Q2").Value = 149.5731411427556
$ws.Range("R2").Value = 1346.1582702848
$ws.Range("S2").Value = 0.110507171133416
$ws.Range("T2").Value = 0.110507171133416

# Row 3 (Target cluster = FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.699506666666667
$ws.Range("H3").Value = 14.09852
$ws.Range("I3").Value = 0.9660495246229048
$ws.Range("J3").Value = 0.9660495246229047
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 85.46317833333335
$ws.Range("N3").Value = 256.389535
$ws.Range("O3").Value = 0.307162904935779
$ws.Range("P3").Value = 0.307162904935779
$ws.Range("Q3").Value = 401.6347763320223
$ws.Range("R3").Value = 3614.712986988201
$ws.Range("S3").Value = 0.2967345782949998
$ws.Range("T3").Value = 0.2967345782949998

# Row 4 (Target cluster = M2)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.699506666666667
$ws.Range("H4").Value = 14.09852
$ws.Range("I4").Value = 0.9660495246229048
$ws.Range("J4").Value = 0.9660495246229047
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 122.2478306666667
$ws.Range("N4").Value = 366.743492
$ws.Range("O4").Value = 0.4393704929064738
$ws.Range("P4").Value = 0.4393704929064738
$ws.Range("Q4").Value = 574.5044952035379
$ws.Range("R4").Value = 5170.54045683184
$ws.Range("S4").Value = 0.4244536558056304
$ws.Range("T4").Value = 0.4244536558056304

# Row 5 (Target cluster = sCs)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.699506666666667
$ws.Range("H5").Value = 14.09852
$ws.Range("I5").Value = 0.9660495246229048
$ws.Range("J5").Value = 0.9660495246229047
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 38.69562533333333
$ws.Range("N5").Value = 116.086876
$ws.Range("O5").Value = 0.1390758092255191
$ws.Range("P5").Value = 0.1390758092255191
$ws.Range("Q5").Value = 181.8503492248356
$ws.Range("R5").Value = 1636.65314302352
$ws.Range("S5").Value = 0.1343541193888586
$ws.Range("T5").Value = 0.1343541193888586

# --- Append new rows 6-9 (Sending cluster = FAPs, Ligand = Sele, Receptor = Cd44) ---
# Row 6 (Target cluster = ECs)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sele"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1651576666666667
$ws.Range("H6").Value = 0.495473
$ws.Range("I6").Value = 0.03395047537709522
$ws.Range("J6").Value = 0.03395047537709522
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 31.82741333333333
$ws.Range("N6").Value = 95.48223999999999
$ws.Range("O6").Value = 0.114390792932228
$ws.Range("P6").Value = 0.114390792932228
$ws.Range("Q6").Value = 5.256541322168888
$ws.Range("R6").Value = 47.30887189951999
$ws.Range("S6").Value = 0.003883621798812005
$ws.Range("T6").Value = 0.003883621798812006

# Row 7 (Target cluster = FAPs)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sele"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1651576666666667
$ws.Range("H7").Value = 0.495473
$ws.Range("I7").Value = 0.03395047537709522
$ws.Range("J7").Value = 0.03395047537709522
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 85.46317833333335
$ws.Range("N7").Value = 256.389535
$ws.Range("O7").Value = 0.307162904935779
$ws.Range("P7").Value = 0.307162904935779
$ws.Range("Q7").Value = 14.11489911945056
$ws.Range("R7").Value = 127.034092075055
$ws.Range("S7").Value = 0.0104283266407792
$ws.Range("T7").Value = 0.0104283266407792

# Row 8 (Target cluster = M2)
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sele"
$ws.Range("C8").Value = "Cd44"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1651576666666667
$ws.Range("H8").Value = 0.495473
$ws.Range("I8").Value = 0.03395047537709522
$ws.Range("J8").Value = 0.03395047537709522
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 122.2478306666667
$ws.Range("N8").Value = 366.743492
$ws.Range("O8").Value = 0.4393704929064738
$ws.Range("P8").Value = 0.4393704929064738
$ws.Range("Q8").Value = 20.19016646796845
$ws.Range("R8").Value = 181.711498211716
$ws.Range("S8").Value = 0.01491683710084343
$ws.Range("T8").Value = 0.01491683710084343

# Row 9 (Target cluster = sCs)
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sele"
$ws.Range("C9").Value = "Cd44"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1651576666666667
$ws.Range("H9").Value = 0.495473
$ws.Range("I9").Value = 0.03395047537709522
$ws.Range("J9").Value = 0.03395047537709522
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 38.69562533333333
$ws.Range("N9").Value = 116.086876
$ws.Range("O9").Value = 0.1390758092255191
$ws.Range("P9").Value = 0.1390758092255191
$ws.Range("Q9").Value = 6.390879190260889
$ws.Range("R9").Value = 57.51791271234799
$ws.Range("S9").Value = 0.004721689836660579
$ws.Range("T9").Value = 0.004721689836660579
